# Edit script for LUMA SizingGuide.xlsx - "Bulldozer CSF" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bulldozer CSF")

# --- Update raw input cells (A,B,C) for rows 2-8 (rows 9-16 left unchanged) ---
# Row 2: A=2, B=2, C=1
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 1

# Row 3: A=4, B=2, C=1
$ws.Range("A3").Value = 4
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1

# Row 4: A=8, B=2, C=1
$ws.Range("A4").Value = 8
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 1

# Row 5: A unchanged (8), B unchanged (4), C=1
$ws.Range("C5").Value = 1

# Row 6: A=8, B=8, C=1
$ws.Range("A6").Value = 8
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = 1

# Row 7: A=10, B=8, C=1
$ws.Range("A7").Value = 10
$ws.Range("B7").Value = 8
$ws.Range("C7").Value = 1

# Row 8: A=12 (B,C unchanged: 8,1)
$ws.Range("A8").Value = 12

# --- Update the "D" label + constants block ---
# A20 changes from "Points per Face" label to new label "D"
$ws.Range("A20").Value = "D"

# B20: 10 -> 8
$ws.Range("B20").Value = 8

# B24: 8 -> 50
$ws.Range("B24").Value = 50

# B25: 4 -> 50
$ws.Range("B25").Value = 50

# B26: 6 -> 1
$ws.Range("B26").Value = 1

# --- Update selection to G19 ---
$ws.Range("G19").Select()

$wb.Save()
